$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D12").Value = "텐서 플로우 스터디 카카오톡 오픈 채팅방 개설 안내"
$ws.Range("E12").Value = "https://tensorflow.blog/2021/03/10/%ed%85%90%ec%84%9c-%ed%94%8c%eb%a1%9c%ec%9a%b0-%ec%8a%a4%ed%84%b0%eb%94%94-%ec%b9%b4%ec%b9%b4%ec%98%a4%ed%86%a1-%ec%98%a4%ed%94%84-%ec%b1%84%ed%8c%85%eb%b0%a9-%ea%b0%9c%ec%84%a4-%ec%95%88%eb%82%b4/"

$ws.Range("D28").Value = "학술논문작성법 강의 소개/정리 요약"
$ws.Range("E28").Value = "https://ropiens.tistory.com/95"

$ws.Range("D39").Value = "Probability concepts explained: Rules of probability"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/Probability-concepts-explained-Rules-of-probability"

$ws.Range("D51").Value = "[세이버메트릭스] 수비 요소를 배제하고 투수의 실력을 평가해보자, DIPS와 FIP"
$ws.Range("E51").Value = "https://bskyvision.com/1142"
